# Updates the cryptos list "Price" (D) and "Volume(1h)" (E) columns on the
# active worksheet for rows 2-51, matching the upstream data refresh.
#
# Many "Price" values (e.g. "0.9973") look like plain numbers to Excel's
# smart-parsing, but the source data stores every value as literal text
# (t="inlineStr"). To keep these cells text (not auto-converted to
# numbers) without leaving a stray number-format/style behind, we briefly
# force a Text format, assign the value, then reset the cell style back to
# "Normal" so the on-disk style index matches the original (unstyled) cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.733.54"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "1.725.40"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9973"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9979"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4842"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2583"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06192"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").Value = "1.725.10"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "15.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06896"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6070"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9980"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "26.543.55"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007159"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").Value = "1.949.44"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.426"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.575"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.068"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.771"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.377"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.974"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07933"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.687"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.597"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6196"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9235"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.022"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.438"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9976"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01496"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.638"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3829"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.851"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1154"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.879"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.80%  "
